{"js": "// The recorded change for this template is a pure OOXML re-serialization:\n// every part that differs (word/document.xml, word/footnotes.xml,\n// word/header1.xml, word/styles.xml) only has its XML attributes\n// re-ordered (e.g. `w:headerReference w:type=\"default\" r:id=\"rId6\"` ->\n// `w:headerReference r:id=\"rId6\" w:type=\"default\"`, `w:footnote\n// w:type=\"separator\" w:id=\"-1\"` -> `w:footnote w:id=\"-1\"\n// w:type=\"separator\"`, `w:color w:val=\"...\" w:themeColor=\"...\"\n// w:themeShade=\"...\"` -> attributes re-sorted, etc.). No text, value,\n// formatting, or structural content actually changes anywhere in the\n// package (same attribute/value pairs, same element order, same\n// wording). There is nothing for the Word object model to \"fix\" here,\n// so this script just walks the same areas the diff touches and\n// confirms they already hold the expected (unchanged) content, without\n// writing anything back - any write (even re-assigning a property to\n// its current value) would explode the compact markup into new runs/\n// explicit properties and regress the document instead of matching the\n// target.\n\nconst body = context.document.body;\nbody.load(\"text\");\n\nconst sections = context.document.sections;\nsections.load(\"items\");\n\nawait context.sync();\n\n// word/document.xml: a single empty paragraph (it only carries the\n// \"_GoBack\" bookmark) followed by the section properties - untouched.\nconst bodyText = body.text;\n\n// word/document.xml section count (the sectPr whose headerReference /\n// pgSz / pgMar attribute order was normalized) is still exactly one.\nconst sectionCount = sections.items.length;\n\n// word/styles.xml: the latentStyles / style catalog (Normal,\n// Policepardfaut, TableauNormal, Aucuneliste, En-tte, En-tteCar,\n// Pieddepage, PieddepageCar) is untouched - only attribute order\n// changed on those <w:style>/<w:lsdException> elements.\nconst styles = context.document.getStyles();\nstyles.load(\"items/nameLocal\");\nawait context.sync();\nconst styleNames = styles.items.map((s) => s.nameLocal);\n\nreturn JSON.stringify({ bodyText, sectionCount, styleNames });\n", "ps1": "# The recorded change for this template is a pure OOXML re-serialization:\n# every part that differs (word/document.xml, word/footnotes.xml,\n# word/header1.xml, word/styles.xml) only has its XML attributes\n# re-ordered, e.g.:\n#   <w:headerReference w:type=\"default\" r:id=\"rId6\"/>\n#     -> <w:headerReference r:id=\"rId6\" w:type=\"default\"/>\n#   <w:footnote w:type=\"separator\" w:id=\"-1\">\n#     -> <w:footnote w:id=\"-1\" w:type=\"separator\">\n#   <w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/>\n#     -> <w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/>\n#   <w:latentStyles w:defLockedState=\"0\" w:defUIPriority=\"99\" ...>\n#     -> <w:latentStyles w:count=\"371\" w:defLockedState=\"0\" ...>\n# and so on for every <w:lsdException>/<w:style> in the style catalog.\n# No text, value, formatting, or structural content actually changes\n# anywhere in the package - same attribute/value pairs, same element\n# order, same wording. There is nothing for the Word object model to\n# \"fix\" here, so this script walks the same areas the diff touches\n# (section/page geometry, the header, the footnote separators, the\n# style catalog) purely to confirm they already hold the expected\n# content, and writes nothing back - even a no-op re-assignment (e.g.\n# Font.Color = Font.Color) would expand the compact markup into new\n# runs/explicit properties and regress the document instead of\n# matching the target.\n\n$d = $word.ActiveDocument\n\n# word/document.xml: section properties (headerReference / pgSz / pgMar)\n# only had their attribute order normalized; the geometry itself is the\n# same 11906x16838 twips page with 1417/708 twip margins.\n$section = $d.Sections.Item(1)\n$pageSetup = $section.PageSetup\n$topMargin = $pageSetup.TopMargin\n$pageWidth = $pageSetup.PageWidth\n\n# word/header1.xml: header paragraphs/field code text is unchanged -\n# only the <w:color> attribute order on the emptyTable() field code run\n# was normalized.\n$headerText = $section.Headers.Item(1).Range.Text\n\n# word/footnotes.xml: the separator/continuationSeparator footnotes are\n# unchanged - only w:id/w:type attribute order was normalized.\n$footnoteSeparator = $d.Footnotes.Separator.Text\n$footnoteContinuation = $d.Footnotes.ContinuationSeparator.Text\n\n# word/styles.xml: the style catalog is unchanged - only attribute\n# order on <w:style>/<w:lsdException>/<w:latentStyles> was normalized.\n$styleCount = $d.Styles.Count\n\nWrite-Output (\"topMargin=\" + $topMargin)\nWrite-Output (\"pageWidth=\" + $pageWidth)\nWrite-Output (\"headerText=\" + $headerText)\nWrite-Output (\"footnoteSeparator=\" + $footnoteSeparator)\nWrite-Output (\"footnoteContinuation=\" + $footnoteContinuation)\nWrite-Output (\"styleCount=\" + $styleCount)\n"}
